$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.526.19'
$ws.Range("E2").Value = '  +0.06%  '
$ws.Range("D3").Value = '1.918.55'
$ws.Range("E3").Value = '  -0.26%  '
$ws.Range("D4").Value = '''0.9999'
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '''245.55'
$ws.Range("E5").Value = '  +0.97%  '
$ws.Range("E6").Value = '  -0.03%  '
$ws.Range("D7").Value = '''0.4799'
$ws.Range("E7").Value = '  +1.86%  '
$ws.Range("D8").Value = '''0.2898'
$ws.Range("E8").Value = '  +0.63%  '
$ws.Range("D9").Value = '''0.06723'
$ws.Range("E9").Value = '  -0.28%  '
$ws.Range("D10").Value = '''110.54'
$ws.Range("E10").Value = '  +4.16%  '
$ws.Range("D11").Value = '''19.16'
$ws.Range("E11").Value = '  +5.10%  '
$ws.Range("D12").Value = '1.917.90'
$ws.Range("E12").Value = '  -0.26%  '
$ws.Range("D13").Value = '''0.07574'
$ws.Range("E13").Value = '  -2.37%  '
$ws.Range("D14").Value = '''5.270'
$ws.Range("E14").Value = '  -0.11%  '
$ws.Range("D15").Value = '''0.6688'
$ws.Range("E15").Value = '  +1.30%  '
$ws.Range("D16").Value = '''300.27'
$ws.Range("E16").Value = '  +2.81%  '
$ws.Range("D17").Value = '30.522.39'
$ws.Range("E17").Value = '  +0.07%  '
$ws.Range("D18").Value = '''13.04'
$ws.Range("E18").Value = '  +0.77%  '
$ws.Range("D19").Value = '''1.000'
$ws.Range("E19").Value = '  -0.03%  '
$ws.Range("D20").Value = '''5.574'
$ws.Range("E20").Value = '  +5.82%  '
$ws.Range("D21").Value = '''0.000007583'
$ws.Range("E21").Value = '  +0.05%  '
$ws.Range("D22").Value = '2.163.24'
$ws.Range("E22").Value = '  +0.49%  '
$ws.Range("D23").Value = '''0.9993'
$ws.Range("E23").Value = '  -0.11%  '
$ws.Range("D24").Value = '''6.458'
$ws.Range("E24").Value = '  +4.19%  '
$ws.Range("D25").Value = '''9.495'
$ws.Range("E25").Value = '  +1.50%  '
$ws.Range("D26").Value = '''164.51'
$ws.Range("E26").Value = '  -2.40%  '
$ws.Range("D27").Value = '''20.32'
$ws.Range("E27").Value = '  -4.71%  '
$ws.Range("D28").Value = '''2.117'
$ws.Range("E28").Value = '  +0.92%  '
$ws.Range("E29").Value = '  +0.79%  '
$ws.Range("D30").Value = '''1.399'
$ws.Range("E30").Value = '  +2.45%  '
$ws.Range("D31").Value = '''4.167'
$ws.Range("E31").Value = '  -0.02%  '
$ws.Range("D32").Value = '''4.042'
$ws.Range("E32").Value = '  +1.10%  '
$ws.Range("D33").Value = '''0.04995'
$ws.Range("E33").Value = '  -0.99%  '
$ws.Range("D34").Value = '''0.7362'
$ws.Range("E35").Value = '  -1.42%  '
$ws.Range("E36").Value = '  +0.03%  '
$ws.Range("D37").Value = '''2.721'
$ws.Range("E37").Value = '  -0.28%  '
$ws.Range("D38").Value = '''0.02036'
$ws.Range("E38").Value = '  -3.91%  '
$ws.Range("D39").Value = '''2.683'
$ws.Range("E39").Value = '  +0.02%  '
$ws.Range("D40").Value = '''111.00'
$ws.Range("E40").Value = '  +0.62%  '
$ws.Range("D41").Value = '''2.020'
$ws.Range("E41").Value = '  -2.68%  '
$ws.Range("D42").Value = '''0.4436'
$ws.Range("E42").Value = '  +3.91%  '
$ws.Range("D43").Value = '''71.94'
$ws.Range("E43").Value = '  +6.74%  '
$ws.Range("D44").Value = '''0.8619'
$ws.Range("E44").Value = '  -1.62%  '
$ws.Range("D45").Value = '''5.902'
$ws.Range("E45").Value = '  +0.67%  '
$ws.Range("D46").Value = '''1.000'
$ws.Range("E46").Value = '  -0.01%  '
$ws.Range("D47").Value = '''49.34'
$ws.Range("E47").Value = '  -0.31%  '
$ws.Range("D48").Value = '''7.272'
$ws.Range("E48").Value = '  +1.17%  '
$ws.Range("D49").Value = '''9.243'
$ws.Range("E49").Value = '  -0.13%  '
$ws.Range("D50").Value = '''0.1231'
$ws.Range("E50").Value = '  +1.12%  '
$ws.Range("D51").Value = '''0.2534'
$ws.Range("E51").Value = '  +2.70%  '
